$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns are stored as text (not numbers/percentages),
# so force each target cell to Text format before assigning the new value.
# This prevents Excel from auto-converting "257.58" -> number or "0.15%" -> percentage.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "257.58"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "0.15%"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "26.75"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "-1.05%"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "4.634"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "0.57%"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.05941"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "0.76%"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.615"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "-0.38%"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.8568"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "-0.91%"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "-1.95%"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1385"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "-1.26%"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.04417"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "15.29%"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07005"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-1.08%"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03019"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "-5.58%"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.09115"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-1.51%"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001521"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "-1.26%"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0006035"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "0.37%"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006107"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "0.27%"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.473"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "-1.22%"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.130"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-1.76%"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "2.151"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "-2.76%"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.3102"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "0.11%"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "1.62%"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.861"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "0.19%"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04182"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-1.16%"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001214"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-0.12%"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004498"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "5.08%"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001197"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "-0.08%"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0001711"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "-11.54%"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.03815"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "-0.36%"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1106"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "0.64%"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.003712"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "-39.95%"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002302"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-0.64%"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.01509"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "29.87%"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005082"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-6.68%"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00000000748"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-0.09%"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.04989"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-35.79%"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "10,478.60%"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00002096"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-0.09%"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0001996"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.09%"
